$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 2, pushing existing rows 2-4 down to 3-5.
$ws.Rows.Item(2).Insert()

# Populate the new row 2 with the machine-readable "slug" identifiers that
# correspond to each column header in row 1, allowing two columns to be
# related to each other to build hierarchical SKOS concepts.
$ws.Cells.Item(2, 1).Value = "tipo-hogar-1"
$ws.Cells.Item(2, 2).Value = "tipo-de-hogar-2"
$ws.Cells.Item(2, 3).Value = "numero-hogares"
$ws.Cells.Item(2, 4).Value = "provincia-codigo"
$ws.Cells.Item(2, 5).Value = "aragon"
$ws.Cells.Item(2, 6).Value = "municipio-codigo"
$ws.Cells.Item(2, 7).Value = "provincia-nombre"
$ws.Cells.Item(2, 8).Value = "municipio-nombre"
